$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E stay as Text so numeric-looking / percent strings are preserved verbatim
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.511.43"

# Row 3
$ws.Range("D3").Value = "2.626.66"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "535.18"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").Value = "142.93"
$ws.Range("E6").Value = "  +1.31%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.23%  "

# Row 9
$ws.Range("D9").Value = "6.98"
$ws.Range("E9").Value = "  +7.94%  "

# Row 10
$ws.Range("E10").Value = "  -1.58%  "

# Row 11
$ws.Range("E11").Value = "  +0.10%  "

# Row 12
$ws.Range("D12").Value = "0.135"
$ws.Range("E12").Value = "  +0.96%  "

# Row 13
$ws.Range("D13").Value = "3.094.26"
$ws.Range("E13").Value = "  +1.14%  "

# Row 14
$ws.Range("D14").Value = "58.461.06"
$ws.Range("E14").Value = "  -1.37%  "

# Row 15
$ws.Range("D15").Value = "20.79"
$ws.Range("E15").Value = "  +0.59%  "

# Row 16
$ws.Range("D16").Value = "2.622.65"
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("E17").Value = "  -0.86%  "

# Row 18
$ws.Range("E18").Value = "  +0.79%  "

# Row 19
$ws.Range("D19").Value = "334.77"
$ws.Range("E19").Value = "  -1.74%  "

# Row 20
$ws.Range("E20").Value = "  +0.69%  "

# Row 21
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("D23").Value = "66.08"
$ws.Range("E23").Value = "  -2.13%  "

# Row 24
$ws.Range("E24").Value = "  +1.78%  "

# Row 25
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.38%  "

# Row 26
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  -1.00%  "

# Row 27
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  -1.04%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0737"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("E30").Value = "  -1.00%  "

# Row 31
$ws.Range("E31").Value = "  +1.02%  "

# Row 32
$ws.Range("D32").Value = "18.75"
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("D33").Value = "150.51"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("E34").Value = "  -1.55%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "37.16"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  +2.37%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.41"
$ws.Range("E38").Value = "  -3.06%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "0.811"
$ws.Range("E39").Value = "  -1.33%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  +1.49%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "280.98"
$ws.Range("E41").Value = "  +3.07%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.08%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.599"
$ws.Range("E43").Value = "  +0.38%  "

# Row 44
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "10.68"
$ws.Range("E44").Value = "  -0.66%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0531"
$ws.Range("E45").Value = "  +1.70%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "18.99"
$ws.Range("E46").Value = "  +2.92%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0935"
$ws.Range("E47").Value = "  -1.80%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0224"
$ws.Range("E48").Value = "  +0.91%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.947.26"
$ws.Range("E49").Value = "  +0.31%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.46"
$ws.Range("E50").Value = "  -0.56%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "17.88"
$ws.Range("E51").Value = "  -3.66%  "

# Restore default (Normal) style on D:E so no residual text-format style lingers on cells
$ws.Range("D2:E51").Style = "Normal"
